$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Result"
$ws.Range("D2").Value = "Failed"
$ws.Range("D3").Borders.LineStyle = 0

$ws.Range("C2").Select()
